$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9444109797477722
$ws.Range("B1").Value = 1.951615214347839
$ws.Range("C1").Value = 4.156449794769287
$ws.Range("D1").Value = 3.277727127075195
$ws.Range("E1").Value = 1.445918798446655
